# Generate Report for Handback
# - The "Ready for handoff" status (Overview!B3/C3, zh-cn!C3, de-de!C3) is
#   now reported as "Handback transform failed" because the handback step
#   failed validation.
# - Each locale sheet gets a new "Error Detail" entry (column K, row 3)
#   explaining the mismatch between the handback file name and the
#   original handoff file name.

$wb = $excel.ActiveWorkbook

$failedStatus = "Handback transform failed"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $failedStatus
$overview.Range("C3").Value = $failedStatus

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $failedStatus
$zhcn.Range("K3").Value = "Handback file name: j2ce5qer.bwb is different with handoff file name: 7e59f1fe-64dc-421a-b255-b4bb6ea16fe0.f2e52a5928e26bd8c86795a2e07631092842f429.zh-cn."

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $failedStatus
$dede.Range("K3").Value = "Handback file name: j2ce5qer.bwb is different with handoff file name: 7e59f1fe-64dc-421a-b255-b4bb6ea16fe0.f2e52a5928e26bd8c86795a2e07631092842f429.de-de."
